# Applies the diff to the Candidates worksheet:
#  - Row 2: candidate fully replaced with new info (Jeenusha John), several
#           screening-related fields cleared out.
#  - Row 3: only Round 1 Remarks (AA3), Round 2 Remarks (AC3) and Screened By
#           (AG3) change.
#  - Row 4: candidate fully replaced with new info (Jeenusha John), application
#           accepted, several previously "-" placeholder fields cleared out.
#
# Some of the new values are purely numeric-looking strings (phone numbers,
# years of experience, CTC figures, notice period) that must stay stored as
# TEXT (the sheet uses inlineStr / text cells throughout), so those cells are
# pre-formatted as Text ("@") before the value is written to stop Excel from
# auto-converting them to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
}

# ---------------------------------------------------------------------------
# Row 2
# ---------------------------------------------------------------------------
$ws.Range("A2").Value  = "2025-12-23T00:00:00.000Z"
$ws.Range("B2").Value  = "Jeenusha John"
$ws.Range("C2").Value  = "jkhjjjqq@example.com"
Set-TextValue $ws.Range("D2") "9999999992"
$ws.Range("E2").Value  = "https://linkedin.com/in/jeenusha"
$ws.Range("F2").Value  = "4a4de728-6ddd-4f48-9441-b3ac101b2291-JeenushaJohn_Resume.pdf"
$ws.Range("G2").Value  = "SRE"
$ws.Range("H2").Value  = "Engineer"
$ws.Range("I2").Value  = "ABC Corp"
Set-TextValue $ws.Range("J2") "5"
$ws.Range("K2").Value  = "Chennai"
$ws.Range("L2").Value  = "Any"
Set-TextValue $ws.Range("M2") "10"
Set-TextValue $ws.Range("N2") "15"
Set-TextValue $ws.Range("O2") "30"
$ws.Range("P2").Value  = "No"
$ws.Range("T2").Value  = ""
$ws.Range("U2").Value  = ""
$ws.Range("V2").Value  = ""
$ws.Range("X2").Value  = ""
$ws.Range("Y2").Value  = ""
$ws.Range("AA2").Value = ""
$ws.Range("AC2").Value = ""

# ---------------------------------------------------------------------------
# Row 3
# ---------------------------------------------------------------------------
$ws.Range("AA3").Value = '{"Communication":"","Technical Assessment":"","Problem-Solving":"","Overall Potential":"","Recommendation":"Proceed Round 2"}'
$ws.Range("AC3").Value = '{"Communication":"","Technical Assessment":"","Problem-Solving":"","Overall Potential":"","Recommendation":" very good ","CTC":""}'
$ws.Range("AG3").Value = "admin"

# ---------------------------------------------------------------------------
# Row 4
# ---------------------------------------------------------------------------
$ws.Range("A4").Value  = "2025-12-23T00:00:00.000Z"
$ws.Range("B4").Value  = "Jeenusha John"
$ws.Range("C4").Value  = "jkhjjj@example.com"
Set-TextValue $ws.Range("D4") "9999999999"
$ws.Range("E4").Value  = "https://linkedin.com/in/jeenusha"
$ws.Range("F4").Value  = "7ae92d4c-aecf-4caf-86cb-61d57e5bdcf3-JeenushaJohn_Resume.pdf"
$ws.Range("G4").Value  = "SRE"
$ws.Range("H4").Value  = "Engineer"
$ws.Range("I4").Value  = "ABC Corp"
Set-TextValue $ws.Range("J4") "5"
$ws.Range("L4").Value  = "Any"
Set-TextValue $ws.Range("M4") "10"
Set-TextValue $ws.Range("N4") "15"
Set-TextValue $ws.Range("O4") "30"
$ws.Range("P4").Value  = "No"
$ws.Range("Q4").Value  = "No"
$ws.Range("T4").Value  = ""
$ws.Range("U4").Value  = ""
$ws.Range("V4").Value  = ""
$ws.Range("W4").Value  = "Accepted"
$ws.Range("X4").Value  = ""
$ws.Range("Y4").Value  = ""
$ws.Range("Z4").Value  = ""
$ws.Range("AB4").Value = ""
$ws.Range("AD4").Value = ""
$ws.Range("AE4").Value = ""
$ws.Range("AF4").Value = "No"
